# daily auto push: 2026-02-01 09:42 UTC
# Insert a new observation row for 2026/02/01 (16:00 bucket) ahead of the
# 2026/12/29 block, pushing every subsequent row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: insert a blank row above the current row 762 (2026/12/29 火).
# This shifts the old rows 762-803 down to 763-804 and grows the sheet
# dimension from D803 to D804 automatically.
$ws.Rows.Item(762).Insert()

# Populate the new row 762. Copy the date cell from the row above (A761,
# which already holds the literal text "2026/02/01") instead of assigning
# the string directly, so Excel doesn't reinterpret it as a date serial.
$ws.Range("A761").Copy()
$ws.Range("A762").PasteSpecial()

$ws.Range("B762").Value = "日"
$ws.Range("C762").Value = 16
$ws.Range("D762").Value = 158
